# Update Name of Algo
# Applies updated RandomForest-imputed values to specific cells on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = -7.446700000000002
$ws.Range("B9").Value = 6.235499999999996
$ws.Range("D12").Value = -7.171899999999998
$ws.Range("D14").Value = -7.792900000000002
$ws.Range("B18").Value = 6.742800000000001
$ws.Range("B20").Value = 8.933499999999995
$ws.Range("D26").Value = -8.368000000000002
$ws.Range("B27").Value = 5.918900000000003
$ws.Range("D27").Value = -8.752599999999999
$ws.Range("D29").Value = -7.384499999999998
$ws.Range("B35").Value = 8.7027
$ws.Range("D37").Value = -7.768799999999995
$ws.Range("D38").Value = -8.0763
$ws.Range("D51").Value = -8.108499999999998
$ws.Range("D52").Value = -7.568300000000001
$ws.Range("D55").Value = -8.676699999999997
$ws.Range("B69").Value = 5.537399999999996
$ws.Range("D69").Value = -7.248899999999999
$ws.Range("D70").Value = -7.353599999999999
$ws.Range("B76").Value = 5.120300000000003
$ws.Range("B78").Value = 9.814500000000004
$ws.Range("D81").Value = -7.673000000000002
$ws.Range("B82").Value = 6.022599999999999
$ws.Range("B83").Value = 5.1764
$ws.Range("D83").Value = -8.894599999999992
$ws.Range("B93").Value = 6.133199999999999
$ws.Range("D102").Value = -7.820999999999998
